# ---------------------------------------------------------------------------
# Stock report correction: quantities (F), computed stock value (G = Rate * Qty),
# a handful of rows whose two lot/batch entries had their Code/MRP/Qty/Value
# fields transposed, and the knock-on "Sub Total:" / "Grand Total:" rows.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity (F) / Value (G) corrections for individual stock lines ---
$ws.Range("F6").Value = 114
$ws.Range("G6").Value = 3406.32
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 384.41
$ws.Range("F82").Value = 4
$ws.Range("G82").Value = 203.2
$ws.Range("F144").Value = 1588
$ws.Range("G144").Value = 13418.6
$ws.Range("F145").Value = 836
$ws.Range("G145").Value = 6679.64
$ws.Range("F196").Value = 1
$ws.Range("G196").Value = 114.56
$ws.Range("F214").Value = 51
$ws.Range("G214").Value = 4472.7
$ws.Range("F248").Value = 0
$ws.Range("G248").Value = 0
$ws.Range("F255").Value = 618
$ws.Range("G255").Value = 105881.94
$ws.Range("F291").Value = 130
$ws.Range("G291").Value = 5591.3
$ws.Range("F328").Value = 68
$ws.Range("G328").Value = 2530.28
$ws.Range("F341").Value = 2
$ws.Range("G341").Value = 101.9
$ws.Range("F434").Value = 47
$ws.Range("G434").Value = 1534.08
$ws.Range("F450").Value = 15
$ws.Range("G450").Value = 2081.1
$ws.Range("F454").Value = 61
$ws.Range("G454").Value = 2083.15
$ws.Range("F492").Value = 70
$ws.Range("G492").Value = 9138.5
$ws.Range("F554").Value = 17
$ws.Range("G554").Value = 633.76
$ws.Range("F580").Value = 73
$ws.Range("G580").Value = 4160.27
$ws.Range("F581").Value = 27
$ws.Range("G581").Value = 6528.6
$ws.Range("F599").Value = 2258
$ws.Range("G599").Value = 368302.38
$ws.Range("F612").Value = 37
$ws.Range("G612").Value = 1516.63

# --- Row pairs whose Code (B) / MRP (E) / Qty (F) / Value (G) were swapped ---
$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 2
$ws.Range("G127").Value = 241.38
$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86
$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72
$ws.Range("B364").Value = 65068
$ws.Range("E364").Value = 13.97
$ws.Range("F364").Value = 63
$ws.Range("G364").Value = 828.45
$ws.Range("B365").Value = 53602
$ws.Range("E365").Value = 15.69
$ws.Range("F365").Value = -231
$ws.Range("G365").Value = -3037.65
$ws.Range("B382").Value = 64919
$ws.Range("E382").Value = 27.97
$ws.Range("F382").Value = 61
$ws.Range("G382").Value = 1604.3
$ws.Range("B383").Value = 45702
$ws.Range("E383").Value = 31.43
$ws.Range("F383").Value = -215
$ws.Range("G383").Value = -5654.5
$ws.Range("B385").Value = 53595
$ws.Range("E385").Value = 17.61
$ws.Range("F385").Value = -335
$ws.Range("G385").Value = -4934.55
$ws.Range("B386").Value = 65067
$ws.Range("E386").Value = 15.65
$ws.Range("F386").Value = 126
$ws.Range("G386").Value = 1855.98
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 5
$ws.Range("G442").Value = 1369.6
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52

# --- Sub Total / Grand Total (B) recalculations ---
$ws.Range("B10").Value = 33443.41
$ws.Range("B90").Value = 209506.96
$ws.Range("B147").Value = 23213.27
$ws.Range("B216").Value = 56402.18
$ws.Range("B260").Value = 216495.33
$ws.Range("B304").Value = 203702.6
$ws.Range("B330").Value = 32996.16
$ws.Range("B346").Value = 30213.01
$ws.Range("B435").Value = 2073.76
$ws.Range("B460").Value = 16467.89
$ws.Range("B493").Value = 15428.56
$ws.Range("B560").Value = 17485.67
$ws.Range("B583").Value = 31558.9
$ws.Range("B606").Value = 562125.84
$ws.Range("B618").Value = 49337.73
$ws.Range("B619").Value = 2101924.11
$ws.Range("B620").Value = 2101924.11
